$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Added price for BLE121: set the qty-1000 unit price for the BLE121 row (row 20).
$ws.Range("I20").Value = 8.3123199999999997

$ws.Range("H23:J23").ClearContents()

# The cell had a yellow "missing price" highlight; clear it now that the
# price has been filled in.
$ws.Range("I20:J20").ClearFormats()
$ws.Range("I20:J20").NumberFormat = '"$"#,##0.000000'

# Update the current selection to reflect where the user left off.
$ws.Range("J23").Select()
